$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug status (row 6, column B) moves from "提交" (submitted) to "分配" (assigned)
$ws.Range("B6").Value = "分配"

# Active cell / selection moves from C14 to C7
$ws.Range("C7").Select()
